# Update "想去人数" (column F) figures on the 展览 and 全部类型 sheets to
# reflect the newly scraped counts (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Values shared by both the "展览" (Exhibition) sheet and the "全部类型"
# (All types) sheet.
$commonUpdates = @{
    "F3"  = 979
    "F4"  = 296
    "F6"  = 1125
    "F8"  = 2430
    "F9"  = 7949
    "F10" = 944
    "F11" = 482
    "F12" = 423
    "F13" = 189
    "F14" = 448
    "F16" = 171
    "F18" = 329
    "F19" = 1416
    "F23" = 198
    "F24" = 348
    "F25" = 192
    "F28" = 118
    "F29" = 37
    "F30" = 435
    "F31" = 1171
    "F32" = 38
    "F33" = 61
    "F35" = 72
    "F36" = 92
    "F38" = 83
    "F39" = 74
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($cellRef in $commonUpdates.Keys) {
        $ws.Range($cellRef).Value = $commonUpdates[$cellRef]
    }
}

# F17 diverged slightly between the two sheets in this refresh.
$wb.Worksheets.Item("展览").Range("F17").Value = 8205
$wb.Worksheets.Item("全部类型").Range("F17").Value = 8206
